$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 248
$ws.Range("H104").Value = 394.2
$ws.Range("H107").Value = 464.93332
$ws.Range("H112").Value = 2224.2856
$ws.Range("H131").Value = 3615.8
$ws.Range("H31").Value = 224.25
$ws.Range("H92").Value = 375.08694
$ws.Range("H99").Value = 420.1875
$ws.Range("I101").Value = 281.83334
$ws.Range("I104").Value = 243
$ws.Range("I107").Value = 339.16666
$ws.Range("I131").Value = 3615.8
$ws.Range("I31").Value = 224.25
$ws.Range("I92").Value = 341.35
$ws.Range("I99").Value = 328.69232
$ws.Range("J112").Value = 2498.4
$ws.Range("J131").Value = 0
$ws.Range("J92").Value = 600
$ws.Range("K101").Value = 845.5000200000001
$ws.Range("K104").Value = 729
$ws.Range("K107").Value = 339.16666
$ws.Range("K131").Value = 10847.4
$ws.Range("K31").Value = 672.75
$ws.Range("K92").Value = 341.35
$ws.Range("K99").Value = 986.07696
$ws.Range("L112").Value = 7495.200000000001
$ws.Range("L131").Value = 0
$ws.Range("L92").Value = 600
$ws.Range("M101").Value = 776.4999799999999
$ws.Range("M104").Value = 1018
$ws.Range("M107").Value = 1580.83334
$ws.Range("M131").Value = -5807.400000000001
$ws.Range("M31").Value = -442.75
$ws.Range("M92").Value = 906.65
$ws.Range("M99").Value = 511.92304
$ws.Range("N112").Value = -9711.200000000001
$ws.Range("N131").ClearContents()
$ws.Range("N92").Value = -3096

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4493.625
$ws.Range("H132").Value = 3699.4546
$ws.Range("H136").Value = 3193.0417
$ws.Range("H32").Value = 7337.136
$ws.Range("H5").Value = 2
$ws.Range("H61").Value = 3193.0417
$ws.Range("H74").Value = 507742
$ws.Range("H77").Value = 507742
$ws.Range("H97").Value = 1603.8889
$ws.Range("I102").Value = 4030.7693
$ws.Range("I132").Value = 3128.2
$ws.Range("I136").Value = 2547.8667
$ws.Range("I32").Value = 6916.684
$ws.Range("I61").Value = 2547.8667
$ws.Range("I97").Value = 1276.4286
$ws.Range("J132").Value = 4175.5
$ws.Range("J5").Value = 2
$ws.Range("J74").Value = 3213.8572
$ws.Range("J77").Value = 3213.8572
$ws.Range("K102").Value = 4030.7693
$ws.Range("K132").Value = 9384.599999999999
$ws.Range("K136").Value = 7643.6001
$ws.Range("K32").Value = 6916.684
$ws.Range("K61").Value = 2547.8667
$ws.Range("K97").Value = 1276.4286
$ws.Range("L132").Value = 12526.5
$ws.Range("L5").Value = 2
$ws.Range("L74").Value = 3213.8572
$ws.Range("L77").Value = 16069.286
$ws.Range("M102").Value = -2408.7693
$ws.Range("M132").Value = -6854.599999999999
$ws.Range("M136").Value = -5093.6001
$ws.Range("M32").Value = -6629.684
$ws.Range("M61").Value = -2335.8667
$ws.Range("M97").Value = -780.4286
$ws.Range("N132").Value = -17586.5
$ws.Range("N5").Value = -226
$ws.Range("N74").Value = -4961.8572
$ws.Range("N77").Value = -24805.286

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13519351
$ws.Range("H4").Value = 2
$ws.Range("I20").Value = 19237426
$ws.Range("J4").Value = 2
$ws.Range("K20").Value = 19237426
$ws.Range("L4").Value = 2
$ws.Range("M20").Value = -19237179
$ws.Range("N4").Value = -232

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1590
$ws.Range("H109").Value = 48861.332
$ws.Range("H132").Value = 5424
$ws.Range("H141").Value = 437263.53
$ws.Range("H19").Value = 809.75
$ws.Range("H22").Value = 2438.5
$ws.Range("H24").Value = 809.75
$ws.Range("H31").Value = 4352.609
$ws.Range("H34").Value = 4352.609
$ws.Range("H50").Value = 68091
$ws.Range("H8").Value = 9999
$ws.Range("H82").Value = 39995
$ws.Range("H85").Value = 39995
$ws.Range("I105").Value = 1036
$ws.Range("I109").Value = 0
$ws.Range("I19").Value = 809.75
$ws.Range("I22").Value = 2526.6
$ws.Range("I24").Value = 809.75
$ws.Range("I31").Value = 3066.4443
$ws.Range("I34").Value = 3066.4443
$ws.Range("I8").Value = 9999
$ws.Range("J109").Value = 48861.332
$ws.Range("J132").Value = 6657.8
$ws.Range("J141").Value = 458496.66
$ws.Range("J31").Value = 8982.799999999999
$ws.Range("J34").Value = 8982.799999999999
$ws.Range("J50").Value = 68091
$ws.Range("J82").Value = 39995
$ws.Range("J85").Value = 39995
$ws.Range("K105").Value = 1036
$ws.Range("K109").Value = 0
$ws.Range("K19").Value = 809.75
$ws.Range("K22").Value = 2526.6
$ws.Range("K24").Value = 809.75
$ws.Range("K31").Value = 3066.4443
$ws.Range("K34").Value = 3066.4443
$ws.Range("K8").Value = 9999
$ws.Range("L109").Value = 48861.332
$ws.Range("L132").Value = 19973.4
$ws.Range("L141").Value = 458496.66
$ws.Range("L31").Value = 8982.799999999999
$ws.Range("L34").Value = 8982.799999999999
$ws.Range("L50").Value = 68091
$ws.Range("L82").Value = 39995
$ws.Range("L85").Value = 39995
$ws.Range("M105").Value = 711
$ws.Range("M109").ClearContents()
$ws.Range("M19").Value = -639.75
$ws.Range("M22").Value = -2176.6
$ws.Range("M24").Value = -639.75
$ws.Range("M31").Value = -2771.4443
$ws.Range("M34").Value = -2864.4443
$ws.Range("M8").Value = -9859
$ws.Range("N109").Value = -50941.332
$ws.Range("N132").Value = -25033.4
$ws.Range("N141").Value = -468856.66
$ws.Range("N31").Value = -9572.799999999999
$ws.Range("N34").Value = -9386.799999999999
$ws.Range("N50").Value = -69341
$ws.Range("N82").Value = -40717
$ws.Range("N85").Value = -42491

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 258
$ws.Range("H2").Value = 1608.3158
$ws.Range("H4").Value = 56703816
$ws.Range("I4").Value = 73520750
$ws.Range("J17").Value = 0
$ws.Range("J2").Value = 3353.2222
$ws.Range("K4").Value = 220562250
$ws.Range("L17").Value = 0
$ws.Range("L2").Value = 20119.3332
$ws.Range("M4").Value = -220562138
$ws.Range("N17").ClearContents()
$ws.Range("N2").Value = -20345.3332

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5799.8
$ws.Range("H70").Value = 82063.62
$ws.Range("H73").Value = 82063.62
$ws.Range("H80").Value = 142861570
$ws.Range("H83").Value = 142861570
$ws.Range("I132").Value = 4999
$ws.Range("I70").Value = 157758.69
$ws.Range("I73").Value = 157758.69
$ws.Range("I80").Value = 250003500
$ws.Range("I83").Value = 250003500
$ws.Range("J132").Value = 6000
$ws.Range("J70").Value = 6368.5386
$ws.Range("J73").Value = 6368.5386
$ws.Range("J80").Value = 5666.3335
$ws.Range("J83").Value = 5666.3335
$ws.Range("K132").Value = 14997
$ws.Range("K70").Value = 157758.69
$ws.Range("K73").Value = 157758.69
$ws.Range("K80").Value = 250003500
$ws.Range("K83").Value = 1250017500
$ws.Range("L132").Value = 18000
$ws.Range("L70").Value = 6368.5386
$ws.Range("L73").Value = 6368.5386
$ws.Range("L80").Value = 5666.3335
$ws.Range("L83").Value = 28331.6675
$ws.Range("M132").Value = -12467
$ws.Range("M70").Value = -157488.69
$ws.Range("M73").Value = -156822.69
$ws.Range("M80").Value = -250002502
$ws.Range("M83").Value = -1250012508
$ws.Range("N132").Value = -23060
$ws.Range("N70").Value = -6908.5386
$ws.Range("N73").Value = -8240.5386
$ws.Range("N80").Value = -7662.3335
$ws.Range("N83").Value = -38315.6675

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2587.6924
$ws.Range("I136").Value = 2393.9312
$ws.Range("J136").Value = 3149.6
$ws.Range("K136").Value = 7181.7936
$ws.Range("L136").Value = 9448.799999999999
$ws.Range("M136").Value = -4631.7936
$ws.Range("N136").Value = -14548.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 55556416
$ws.Range("H136").Value = 7832.74
$ws.Range("H82").Value = 50000
$ws.Range("H85").Value = 50000
$ws.Range("I100").Value = 767.2308
$ws.Range("I136").Value = 7066.237
$ws.Range("J136").Value = 10260
$ws.Range("J82").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("K100").Value = 1534.4616
$ws.Range("K136").Value = 21198.711
$ws.Range("L136").Value = 30780
$ws.Range("L82").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("M100").Value = -993.4616000000001
$ws.Range("M136").Value = -18648.711
$ws.Range("N136").Value = -35880
$ws.Range("N82").Value = -50766
$ws.Range("N85").Value = -52652

Write-Host "Applied all market-data + profit updates."